# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# ---------------------------------------------------------------------------
# "Bad Drivers" table (rows 3-5): two existing driver rows got new rollup
# numbers this week, and the Totals row needs to reflect their sum.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Intel(R) Wi-Fi 7 BE200 320MHz - 23.110.0.5"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 144
$ws.Range("D3").Value = 96.90000000000001

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.150.0.3"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 450
$ws.Range("D4").Value = 98.90000000000001

$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 594

# ---------------------------------------------------------------------------
# "Good Drivers" table (rows 13-22 -> 13-23): a brand-new driver version
# showed up this week and is now the top row of the table, so every
# existing row shifts down by one. Insert a fresh row at 13 (shifting the
# rest down, row 22's driver lands on the new row 23) and pick up its
# formatting from the row that used to be on top.
# ---------------------------------------------------------------------------
$ws.Rows("13:13").Insert(-4121)

$ws.Range("A14:E14").Copy() | Out-Null
$ws.Range("A13:E13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B13").Value = 11128
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = ""

# Rows 14-22 keep the same drivers as before (now one row lower) but several
# picked up updated sample counts / good-roaming percentages this week.
$ws.Range("B14").Value = 486214
$ws.Range("D14").Value = 99.90000000000001

$ws.Range("B15").Value = 55507
$ws.Range("D15").Value = 100

$ws.Range("D17").Value = 99.90000000000001

$ws.Range("D18").Value = 100

$ws.Range("B19").Value = 79953
$ws.Range("D19").Value = 99.90000000000001

$ws.Range("B20").Value = 35355

$ws.Range("B21").Value = 65425

$ws.Range("B22").Value = 117653
